# Add description/data-driven "Login" worksheet after the existing "Signup" sheet,
# populate it with Email/Password header + sample credential row (with mailto hyperlinks),
# and update the active-sheet/selection view state on both sheets.

$wb = $excel.ActiveWorkbook
$signup = $wb.Worksheets.Item("Signup")

# Create the new "Login" sheet positioned right after "Signup".
# We add a throwaway sheet first and delete it so the surviving "Login" sheet
# receives the next-next internal sheetId (matching a workbook whose sheetId 2
# was used and removed earlier in its edit history).
$placeholder = $wb.Worksheets.Add($null, $signup)
$placeholder.Name = "TempPlaceholder"
$loginNew = $wb.Worksheets.Add($null, $placeholder)
$loginNew.Name = "Login"
$placeholder.Delete()

# Re-fetch a fresh reference to the Login sheet (avoids a stale reference after the delete).
$login = $wb.Worksheets.Item("Login")

# Header row
$login.Range("A1").Value = "Email"
$login.Range("B1").Value = "Password"

# Sample data row
$login.Range("A2").Value = "test@gmail.com"
$login.Range("B2").Value = "Test@12345"

# Column widths (matching Signup's column sizing as closely as this engine allows)
$login.Columns.Item(1).ColumnWidth = 14.6667
$login.Columns.Item(2).ColumnWidth = 12.8333

# Hyperlink the sample data cells (mailto:, mirroring the pattern already used on Signup)
$login.Hyperlinks.Add($login.Range("A2"), "mailto:test@gmail.com")
$login.Hyperlinks.Add($login.Range("B2"), "mailto:Test@12345")

# Re-apply the workbook's existing "Hyperlink" cell format (copied from Signup!C2) so the
# new cells reuse the already-defined Hyperlink style instead of letting the hyperlink
# creation leave any stray ad-hoc formatting behind.
$signup.Range("C2").Copy()
$login.Range("A2:B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update view/selection state: Signup is no longer the active tab and its selection moves to A5.
$signup.Select()
$signup.Range("A5").Select()

# Login becomes the active/selected tab, with B11 selected - applied last so it "wins" as the
# workbook's active sheet.
$login.Select()
$login.Range("B11").Select()
